$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Día 1: horas estimadas (G) y horas consumidas (H) para las dos primeras tareas
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 3
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1

# Selección activa final en la vista
$ws.Range("H10").Select()
